# "Generate Report for handoff"
#
# Renames the in-flight handoff markdown file (fe270f20-... -> 83cd5807-...)
# wherever it is referenced, introduces a second handoff file
# (15c23a97-...md) whose transform failed, and pushes the previously-last
# ".localization-config" row down to make room for it. New zh-cn / de-de
# xlf file names + handoff timestamps are recorded on the per-language
# sheets as well.

$wb = $excel.ActiveWorkbook

$oldFileUuid = "fe270f20-f63a-4c52-8022-b59a8bd76f63"
$newFileUuid = "83cd5807-c5c1-4023-8817-42705f68e534"
$newConfigUuid = "15c23a97-954e-40e2-94d1-474efd7cde30"
$oldHash = "4bc372c7db1fa7e10b33a2d7313212cdebf7941e"
$newHash = "073c6649ac59146b86d291fa5f21771647aaef62"

$newFileName = "$newFileUuid.md"
$newConfigFileName = "$newConfigUuid.md"
$configDisplay = ".localization-config"

$zhXlfName = "$newFileUuid.$newHash.zh-cn.xlf"
$deXlfName = "$newFileUuid.$newHash.de-de.xlf"

$zhHandoffTime = "2016-01-13 02:33:27"
$deHandoffTime = "2016-01-13 02:33:47"

$epoch = "0001-01-01 00:00:00"
$failedStatus = "Handoff transform failed"
$notLocalized = "Not to be localized"

# Same cornflower-blue underline used by the existing file-name hyperlink
# cells (matches the workbook's "HyperLink" cell style / FF6495ED font).
$linkColor = 15570276
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# Base URLs (rewritten to point at the new file names, same repo layout).
$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/8222858a1db92d8a3ee7e72dbd6146f9ca3ce014/e2e"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8222858a1db92d8a3ee7e72dbd6146f9ca3ce014/.localization-config"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0437a477a8486c396e33af851b46bad0dc8d367c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8df3e88e7d3aa402c0a3edb913d17c8c7646d0bf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName"

function Retarget-Hyperlink($ws, $a1, $newAddress, $newDisplay) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address(0, 0) -eq $a1) {
            $hl.Address = $newAddress
            $hl.TextToDisplay = $newDisplay
            return
        }
    }
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 2: rename handoff markdown file.
$ws1.Range("A2").Value = $newFileName
Retarget-Hyperlink $ws1 "A2" "$mdUrlBase/$newFileName" $newFileName

# Row 3 used to be ".localization-config" -- it now becomes the new,
# failed-transform handoff file; ".localization-config" moves to row 4.
$ws1.Range("A3").Value = $newConfigFileName
$ws1.Range("B3").Value = $failedStatus
$ws1.Range("C3").Value = $failedStatus
Retarget-Hyperlink $ws1 "A3" "$mdUrlBase/$newConfigFileName" $newConfigFileName

# Row 4 (new): the ".localization-config" row.
$ws1.Range("A4").Value = $configDisplay
$ws1.Range("B4").Value = $notLocalized
$ws1.Range("C4").Value = $notLocalized
$ws1.Hyperlinks.Add($ws1.Range("A4"), $configUrl, "", "", $configDisplay) | Out-Null
# Hyperlinks.Add stamps its own theme-coloured style on the cell; reapply
# the workbook's existing cornflower-blue "HyperLink" look on top of it.
$ws1.Range("A4").Style = "HyperLink"
$ws1.Range("A4").Font.Underline = $true
$ws1.Range("A4").Font.Color = $linkColor

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newFileName
$ws2.Range("C2").Value = $zhXlfName
$ws2.Range("D2").Value = $zhHandoffTime
Retarget-Hyperlink $ws2 "A2" "$mdUrlBase/$newFileName" $newFileName
Retarget-Hyperlink $ws2 "C2" $zhXlfUrl $zhXlfName

$ws2.Range("A3").Value = $newConfigFileName
$ws2.Range("B3").Value = $failedStatus
Retarget-Hyperlink $ws2 "A3" "$mdUrlBase/$newConfigFileName" $newConfigFileName

$ws2.Range("A4").Value = $configDisplay
$ws2.Range("B4").Value = $notLocalized
$ws2.Range("D4").Value = $epoch
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = "Ignored"
$ws2.Range("D4").NumberFormat = $dateFormat
$ws2.Hyperlinks.Add($ws2.Range("A4"), $configUrl, "", "", $configDisplay) | Out-Null
$ws2.Range("A4").Style = "HyperLink"
$ws2.Range("A4").Font.Underline = $true
$ws2.Range("A4").Font.Color = $linkColor

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newFileName
$ws3.Range("C2").Value = $deXlfName
$ws3.Range("D2").Value = $deHandoffTime
Retarget-Hyperlink $ws3 "A2" "$mdUrlBase/$newFileName" $newFileName
Retarget-Hyperlink $ws3 "C2" $deXlfUrl $deXlfName

$ws3.Range("A3").Value = $newConfigFileName
$ws3.Range("B3").Value = $failedStatus
Retarget-Hyperlink $ws3 "A3" "$mdUrlBase/$newConfigFileName" $newConfigFileName

$ws3.Range("A4").Value = $configDisplay
$ws3.Range("B4").Value = $notLocalized
$ws3.Range("D4").Value = $epoch
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = "Ignored"
$ws3.Range("D4").NumberFormat = $dateFormat
$ws3.Hyperlinks.Add($ws3.Range("A4"), $configUrl, "", "", $configDisplay) | Out-Null
$ws3.Range("A4").Style = "HyperLink"
$ws3.Range("A4").Font.Underline = $true
$ws3.Range("A4").Font.Color = $linkColor
